$d = $word.ActiveDocument

# --- Step 1: remove the "Jatek alljon meg es kerdezze meg el e szeretned
# kezdeni?" confirmation line entirely (text + its paragraph mark). Locate
# it by content rather than a hard-coded index so the script is resilient.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Jatek alljon meg*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the 'Jatek alljon meg...' paragraph"
}

$confirmPara = $d.Paragraphs.Item($targetIndex)
$confirmPara.Range.Delete()

# Deleting the paragraph's Range (text + paragraph mark) merges it away; the
# paragraph that used to trail it (an empty, bold-formatted paragraph) now
# slides up into its place at $targetIndex.
$mergedPara = $d.Paragraphs.Item($targetIndex)

# --- Step 2: make sure that surviving paragraph is a clean, empty, bold
# paragraph (no stray inherited text/pStyle), exactly like a freshly
# authored one.
$cleanParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$mergedPara.Range.InsertXML($cleanParagraphXml) | Out-Null

$bookmarkPara = $d.Paragraphs.Item($targetIndex)

# --- Step 3: relocate the "_GoBack" bookmark (previously sitting right
# after "Jatek"/"Jatek") into this now-empty paragraph. Re-adding a
# bookmark with the same name moves it off its old location.
$d.Bookmarks.Add("_GoBack", $bookmarkPara.Range) | Out-Null
